# Update the Agt-Agtr1a LR-pairs sheet with refreshed TPM-derived values.
#
# The underlying NATMI computation now includes a third cell type ("ECs")
# as both a sending and a target cluster (previously only "FAPs" and
# "MuSCs" appeared), so the 2 (senders) x 3 (targets) = 6 data rows grow
# into a 3 x 3 = 9 data row table (rows 2-10 on Sheet1, below the header
# in row 1). All of the numeric edge-expression statistics are
# recalculated against the new TPM values as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-T for each data row:
# Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate,
# Ligand average expression value, Ligand total expression value,
# Ligand derived specificity of average expression value,
# Ligand derived specificity of total expression value,
# Receptor-expressing cells, Receptor detection rate,
# Receptor average expression value, Receptor total expression value,
# Receptor derived specificity of average expression value,
# Receptor derived specificity of total expression value,
# Edge average expression weight, Edge total expression weight,
# Edge average expression derived specificity,
# Edge total expression derived specificity
$data = @(
  @("ECs",   "Agt", "Agtr1a", "ECs",   2, 0.6666666666666666, 0.1214023333333333, 0.364207, 0.2856182748266287, 0.2856182748266287, 3, 1, 5.922617666666667,  17.767853,  0.2173916203328182, 0.2173916203328182, 0.7190196041745556, 6.471176437571,     0.06209101956122499, 0.06209101956122498),
  @("ECs",   "Agt", "Agtr1a", "FAPs",  2, 0.6666666666666666, 0.1214023333333333, 0.364207, 0.2856182748266287, 0.2856182748266287, 3, 1, 9.221608999999999,  27.664827,  0.3384821772083041, 0.3384821772083041, 1.119524849687667,  10.075723647189,    0.09667669551379703, 0.09667669551379704),
  @("ECs",   "Agt", "Agtr1a", "MuSCs", 2, 0.6666666666666666, 0.1214023333333333, 0.364207, 0.2856182748266287, 0.2856182748266287, 3, 1, 12.099775,          36.299325,  0.4441262024588777, 0.4441262024588777, 1.468940917808333,  13.220468260275,    0.1268505597516066,  0.1268505597516066),
  @("FAPs",  "Agt", "Agtr1a", "ECs",   1, 0.3333333333333333, 0.174539,           0.523617, 0.4106307243130825, 0.4106307243130825, 3, 1, 5.922617666666667,  17.767853,  0.2173916203328182, 0.2173916203328182, 1.033727764922333,  9.303549884300999,  0.08926767851685977, 0.08926767851685975),
  @("FAPs",  "Agt", "Agtr1a", "FAPs",  1, 0.3333333333333333, 0.174539,           0.523617, 0.4106307243130825, 0.4106307243130825, 3, 1, 9.221608999999999,  27.664827,  0.3384821772083041, 0.3384821772083041, 1.609530413251,     14.485773719259,    0.138991181594115,   0.1389911815941151),
  @("FAPs",  "Agt", "Agtr1a", "MuSCs", 1, 0.3333333333333333, 0.174539,           0.523617, 0.4106307243130825, 0.4106307243130825, 3, 1, 12.099775,          36.299325,  0.4441262024588777, 0.4441262024588777, 2.111882628725,     19.006943658525,    0.1823718642021076,  0.1823718642021076),
  @("MuSCs", "Agt", "Agtr1a", "ECs",   3, 1,                  0.1291096666666667, 0.387329, 0.3037510008602889, 0.3037510008602889, 3, 1, 5.922617666666667,  17.767853,  0.2173916203328182, 0.2173916203328182, 0.7646671927374443, 6.882004734636999,  0.06603292225473346, 0.06603292225473345),
  @("MuSCs", "Agt", "Agtr1a", "FAPs",  3, 1,                  0.1291096666666667, 0.387329, 0.3037510008602889, 0.3037510008602889, 3, 1, 9.221608999999999,  27.664827,  0.3384821772083041, 0.3384821772083041, 1.190598864120333,  10.715389777083,    0.102814300100392,   0.1028143001003921),
  @("MuSCs", "Agt", "Agtr1a", "MuSCs", 3, 1,                  0.1291096666666667, 0.387329, 0.3037510008602889, 0.3037510008602889, 3, 1, 12.099775,          36.299325,  0.4441262024588777, 0.4441262024588777, 1.562197916991667,  14.059781252925,    0.1349037785051634,  0.1349037785051634)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
